# Generate Report for Handback
# Simulates a fresh CI run of the handback-status report: new source GUIDs,
# a new merged xlf content-hash, and refreshed timestamps replace the
# previous run's values across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba"
$oldGuid2 = "63878f5d-2025-4b80-aa58-1777943db14a"
$newGuid1 = "0b5f6bce-de4f-4ab0-a357-3bffd6d9522f"
$newGuid2 = "ffff0af203f0-5743-4a77-9f02-36c49ff57f85"

$newHash = "64b4fca6f110a0dca498345e2827d13cdc8870a1"

$newGenDate = "2016-09-04 19:07:39"
$newZhStart = "2016-09-04 19:07:35"
$newZhEnd   = "2016-09-04 19:07:53"
$newDeEnd   = "2016-09-04 19:08:02"

$file1Name = "$newGuid1.md"
$file2Name = "$newGuid2.md"
$file1Disp = "e2e\$newGuid1.md"
$file2Disp = "e2e\$newGuid2.md"

$zhXlf = "$newGuid1.$newHash.zh-cn.xlf"
$deXlf = "$newGuid1.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $file1Name
$wsOverview.Range("B2").Value = $file1Disp
$wsOverview.Range("G2").Value = $newGenDate

$wsOverview.Range("A3").Value = $file2Name
$wsOverview.Range("B3").Value = $file2Disp
$wsOverview.Range("G3").Value = $newGenDate

# Refresh the two hyperlinks on this sheet so their displayed text matches
# the new file names while keeping the exact same link targets.
$rIdUrls = @(
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid1.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid2.md"
)
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $rIdUrls[0], $null, $null, $file1Disp)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $rIdUrls[1], $null, $null, $file2Disp)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $file1Name
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $newZhStart
$wsZh.Range("I2").Value = $file1Name
$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $newZhEnd

$wsZh.Range("A3").Value = $file2Name
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $newZhStart
$wsZh.Range("I3").Value = $file2Name
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $newZhEnd

$zhUrls = @(
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid1.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6cc27c134a01cb8071f1aba242e89725c73faf33/e2e/$oldGuid1.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid2.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6cc27c134a01cb8071f1aba242e89725c73faf33/e2e/$oldGuid2.md"
)
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhUrls[0], $null, $null, $file1Name)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhUrls[1], $null, $null, $file1Name)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhUrls[2], $null, $null, $file2Name)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhUrls[3], $null, $null, $file2Name)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $file1Name
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $newGenDate
$wsDe.Range("I2").Value = $file1Name
$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $newDeEnd

$wsDe.Range("A3").Value = $file2Name
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $newGenDate
$wsDe.Range("I3").Value = $file2Name
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $newDeEnd

$deUrls = @(
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid1.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f6739e9581cf98be7b04498e032ddece8ea674d/e2e/$oldGuid1.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea4910bc8e87e2b148f96ea187c523dd0caad420/e2e/$oldGuid2.md",
  "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f6739e9581cf98be7b04498e032ddece8ea674d/e2e/$oldGuid2.md"
)
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deUrls[0], $null, $null, $file1Name)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deUrls[1], $null, $null, $file1Name)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deUrls[2], $null, $null, $file2Name)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deUrls[3], $null, $null, $file2Name)

$wb.Save()
